$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the 3 new "2021" rows into the yearly excess-mortality table ---
#
# Original layout: Switzerland (rows 2-12), Sweden (rows 13-23), Spain (rows 24-31).
# New layout: each country gets an additional 2021 row, appended at the end of
# its own block, i.e. just before the next country's block begins (or at the
# very end of the sheet for the last country).
#
# Step 1: insert a new row right before Sweden's block (old row 13) to hold
# Switzerland's 2021 figures.
$ws.Rows.Item(13).Insert()

# Step 2: insert a new row right before Spain's block. After step 1, Spain's
# first row (originally row 24) now sits at row 25, so inserting there makes
# room for Sweden's 2021 figures at the end of Sweden's block.
$ws.Rows.Item(25).Insert()

# After the two inserts, Spain's 8 original rows (originally 24-31) now occupy
# rows 26-33, so row 34 is the first free row after Spain's block - that's
# where Spain's new 2021 figures go (no insert needed, just fill the row).

# --- Row 13: Switzerland, 2021 ---
$ws.Range("A13").Value = "Switzerland"
$ws.Range("B13").Value = 2021
$ws.Range("C13").Value = -803
$ws.Range("D13").Value = -668
$ws.Range("E13").Value = -2696
$ws.Range("F13").Value = 1365
$ws.Range("G13").Value = 125
$ws.Range("H13").Value = -1850
$ws.Range("I13").Value = 2080

# --- Row 25: Sweden, 2021 (only the Global-Serfling columns are populated) ---
$ws.Range("A25").Value = "Sweden"
$ws.Range("B25").Value = 2021
$ws.Range("C25").Value = 1919
$ws.Range("D25").Value = 1081
$ws.Range("E25").Value = -1496
$ws.Range("F25").Value = 3541
# G25 / H25 / I25 intentionally left blank - no Age-Serfling estimate yet for 2021.

# --- Row 34: Spain, 2021 ---
$ws.Range("A34").Value = "Spain"
$ws.Range("B34").Value = 2021
$ws.Range("C34").Value = 11317
$ws.Range("D34").Value = 11860
$ws.Range("E34").Value = -3219
$ws.Range("F34").Value = 26026
$ws.Range("G34").Value = 13163
$ws.Range("H34").Value = -292
$ws.Range("I34").Value = 26482
